$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (linear)
$ws.Range("B4").Value = 0.6464299407097029
$ws.Range("C4").Value = 0.6690000000000002
$ws.Range("D4").Value = 0.6726204350704883
$ws.Range("E4").Value = 0.65
$ws.Range("J4").Value = 0.6814720846263654
$ws.Range("K4").Value = 0.6970000000000001
$ws.Range("L4").Value = 0.6872220031814477
$ws.Range("M4").Value = 0.6845

# Row 5 (rbf)
$ws.Range("B5").Value = 0.6638522343565456
$ws.Range("C5").Value = 0.655
$ws.Range("D5").Value = 0.720816064029648
$ws.Range("E5").Value = 0.6995
$ws.Range("F5").Value = 0.6495272800975109
$ws.Range("G5").Value = 0.9520000000000002
$ws.Range("H5").Value = 0.4932141918756495
$ws.Range("I5").Value = 0.489
$ws.Range("J5").Value = 0.6732384748072968
$ws.Range("K5").Value = 0.6789999999999999
$ws.Range("L5").Value = 0.6910129273954774
$ws.Range("M5").Value = 0.6805

# Row 6 (poly)
$ws.Range("B6").Value = 0.6286100584114255
$ws.Range("C6").Value = 0.5980000000000001
$ws.Range("D6").Value = 0.7261531982035309
$ws.Range("E6").Value = 0.6849999999999998
$ws.Range("J6").Value = 0.6686910638764928
$ws.Range("K6").Value = 0.66
$ws.Range("L6").Value = 0.7069230812855634
$ws.Range("M6").Value = 0.6855

# Row 7 (sigmoid)
$ws.Range("B7").Value = 0.6294446968967373
$ws.Range("C7").Value = 0.667
$ws.Range("D7").Value = 0.612703039295774
$ws.Range("E7").Value = 0.608
$ws.Range("J7").Value = 0.5667790773836666
$ws.Range("K7").Value = 0.6050000000000001
$ws.Range("L7").Value = 0.5807291844353841
$ws.Range("M7").Value = 0.5730000000000001
